# Scheduled market-price refresh for the Pandaemonium Leve profit tracker.
# Updates the cached Universalis price/profit columns (H:N) per Leve row;
# columns A:G (leve metadata) are untouched.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets("ALC")
$ws.Range("H6").Value = 476
$ws.Range("I6").Value = 345
$ws.Range("J6").Value = 1000
$ws.Range("K6").Value = 1035
$ws.Range("L6").Value = 3000
$ws.Range("M6").Value = -923
$ws.Range("N6").Value = -3224

$ws.Range("H69").Value = 4174.1665
$ws.Range("J69").Value = 3009
$ws.Range("L69").Value = 9027
$ws.Range("N69").Value = -10775

$ws.Range("H70").Value = 1750.2051
$ws.Range("I70").Value = 1621.8695
$ws.Range("J70").Value = 1934.6875
$ws.Range("K70").Value = 4865.6085
$ws.Range("L70").Value = 5804.0625
$ws.Range("M70").Value = -4595.6085
$ws.Range("N70").Value = -6344.0625

$ws.Range("H72").Value = 4174.1665
$ws.Range("J72").Value = 3009
$ws.Range("L72").Value = 27081
$ws.Range("N72").Value = -35817

$ws.Range("H73").Value = 1750.2051
$ws.Range("I73").Value = 1621.8695
$ws.Range("J73").Value = 1934.6875
$ws.Range("K73").Value = 4865.6085
$ws.Range("L73").Value = 5804.0625
$ws.Range("M73").Value = -3929.6085
$ws.Range("N73").Value = -7676.0625

$ws.Range("H75").Value = 29800
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()

$ws.Range("H76").Value = 4042.121
$ws.Range("I76").Value = 3644.074
$ws.Range("K76").Value = 3644.074
$ws.Range("M76").Value = -3329.074

$ws.Range("H78").Value = 29800
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()

$ws.Range("H79").Value = 4042.121
$ws.Range("I79").Value = 3644.074
$ws.Range("K79").Value = 3644.074
$ws.Range("M79").Value = -2552.074

$ws.Range("H112").Value = 6481
$ws.Range("J112").Value = 1816.7273
$ws.Range("L112").Value = 5450.1819
$ws.Range("N112").Value = -7666.1819

$ws.Range("H137").Value = 597783.4399999999
$ws.Range("I137").Value = 2597.68
$ws.Range("J137").Value = 1077772
$ws.Range("K137").Value = 7793.039999999999
$ws.Range("L137").Value = 3233316
$ws.Range("M137").Value = -5243.039999999999
$ws.Range("N137").Value = -3238416

$ws.Range("H138").Value = 4268.9697
$ws.Range("I138").Value = 1911.9333
$ws.Range("J138").Value = 4962.216
$ws.Range("K138").Value = 5735.7999
$ws.Range("L138").Value = 14886.648
$ws.Range("M138").Value = -595.7999
$ws.Range("N138").Value = -25166.648

$ws = $wb.Worksheets("ARM")
$ws.Range("H32").Value = 20915.8
$ws.Range("I32").Value = 22306.62
$ws.Range("J32").Value = 7007.6
$ws.Range("K32").Value = 22306.62
$ws.Range("L32").Value = 7007.6
$ws.Range("M32").Value = -22019.62
$ws.Range("N32").Value = -7581.6

$ws.Range("H61").Value = 8282.392
$ws.Range("I61").Value = 6797.7144
$ws.Range("J61").Value = 10591.889
$ws.Range("K61").Value = 6797.7144
$ws.Range("L61").Value = 10591.889
$ws.Range("M61").Value = -6585.7144
$ws.Range("N61").Value = -11015.889

$ws.Range("H88").Value = 5858.7334
$ws.Range("I88").Value = 14125
$ws.Range("J88").Value = 2852.818
$ws.Range("K88").Value = 14125
$ws.Range("L88").Value = 2852.818
$ws.Range("M88").Value = -13719
$ws.Range("N88").Value = -3664.818

$ws.Range("H91").Value = 5858.7334
$ws.Range("I91").Value = 14125
$ws.Range("J91").Value = 2852.818
$ws.Range("K91").Value = 14125
$ws.Range("L91").Value = 2852.818
$ws.Range("M91").Value = -12721
$ws.Range("N91").Value = -5660.818

$ws.Range("H136").Value = 8282.392
$ws.Range("I136").Value = 6797.7144
$ws.Range("J136").Value = 10591.889
$ws.Range("K136").Value = 20393.1432
$ws.Range("L136").Value = 31775.667
$ws.Range("M136").Value = -17843.1432
$ws.Range("N136").Value = -36875.667

$ws = $wb.Worksheets("BSM")
$ws.Range("H105").Value = 7065.273
$ws.Range("I105").Value = 10977.5
$ws.Range("K105").Value = 10977.5
$ws.Range("M105").Value = -9230.5

$ws.Range("H134").Value = 2537.889
$ws.Range("I134").Value = 2537.889
$ws.Range("K134").Value = 7613.667
$ws.Range("M134").Value = -5078.667

$ws = $wb.Worksheets("CRP")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()

$ws.Range("H31").Value = 732189.75
$ws.Range("I31").Value = 12408.765
$ws.Range("J31").Value = 1092080.2
$ws.Range("K31").Value = 12408.765
$ws.Range("L31").Value = 1092080.2
$ws.Range("M31").Value = -12113.765
$ws.Range("N31").Value = -1092670.2

$ws.Range("H34").Value = 732189.75
$ws.Range("I34").Value = 12408.765
$ws.Range("J34").Value = 1092080.2
$ws.Range("K34").Value = 12408.765
$ws.Range("L34").Value = 1092080.2
$ws.Range("M34").Value = -12206.765
$ws.Range("N34").Value = -1092484.2

$ws.Range("H58").Value = 1979998.4
$ws.Range("I58").Value = 3368253.5
$ws.Range("K58").Value = 3368253.5
$ws.Range("M58").Value = -3368050.5

$ws.Range("H132").Value = 2973.282
$ws.Range("I132").Value = 2607.6875
$ws.Range("K132").Value = 7823.0625
$ws.Range("M132").Value = -5293.0625

$ws.Range("H134").Value = 2121.186
$ws.Range("I134").Value = 1596.4546
$ws.Range("J134").Value = 3852.8
$ws.Range("K134").Value = 4789.3638
$ws.Range("L134").Value = 11558.4
$ws.Range("M134").Value = -2254.3638
$ws.Range("N134").Value = -16628.4

$ws.Range("H136").Value = 1979998.4
$ws.Range("I136").Value = 3368253.5
$ws.Range("K136").Value = 10104760.5
$ws.Range("M136").Value = -10102210.5

$ws = $wb.Worksheets("CUL")
$ws.Range("H3").Value = 3325.65
$ws.Range("I3").Value = 1998.9474
$ws.Range("J3").Value = 4526
$ws.Range("K3").Value = 5996.8422
$ws.Range("L3").Value = 13578
$ws.Range("M3").Value = -5884.8422
$ws.Range("N3").Value = -13802

$ws.Range("H133").Value = 3200.0454
$ws.Range("I133").Value = 2102.3076
$ws.Range("J133").Value = 4785.6665
$ws.Range("K133").Value = 6306.9228
$ws.Range("L133").Value = 14356.9995
$ws.Range("M133").Value = -1246.9228
$ws.Range("N133").Value = -24476.9995

$ws = $wb.Worksheets("GSM")
$ws.Range("H82").Value = 37531.75
$ws.Range("J82").Value = 37531.75
$ws.Range("L82").Value = 37531.75
$ws.Range("N82").Value = -38297.75

$ws.Range("H85").Value = 37531.75
$ws.Range("J85").Value = 37531.75
$ws.Range("L85").Value = 37531.75
$ws.Range("N85").Value = -40183.75

$ws.Range("H132").Value = 2684
$ws.Range("I132").Value = 1885.8334
$ws.Range("J132").Value = 3420.7693
$ws.Range("K132").Value = 5657.5002
$ws.Range("L132").Value = 10262.3079
$ws.Range("M132").Value = -3127.5002
$ws.Range("N132").Value = -15322.3079

$ws = $wb.Worksheets("LTW")
$ws.Range("H40").Value = 3675.625
$ws.Range("I40").Value = 3625
$ws.Range("J40").Value = 3726.25
$ws.Range("K40").Value = 3625
$ws.Range("L40").Value = 3726.25
$ws.Range("M40").Value = -3489
$ws.Range("N40").Value = -3998.25

$ws.Range("H108").Value = 79800
$ws.Range("J108").Value = 79800
$ws.Range("L108").Value = 79800
$ws.Range("N108").Value = -87480

$ws.Range("H132").Value = 5626.41
$ws.Range("I132").Value = 6707.6553
$ws.Range("J132").Value = 2490.8
$ws.Range("K132").Value = 20122.9659
$ws.Range("L132").Value = 7472.400000000001
$ws.Range("M132").Value = -17592.9659
$ws.Range("N132").Value = -12532.4

$ws = $wb.Worksheets("WVR")
$ws.Range("H62").Value = 3953.75
$ws.Range("I62").Value = 3982.5
$ws.Range("J62").Value = 3925
$ws.Range("K62").Value = 3982.5
$ws.Range("L62").Value = 3925
$ws.Range("M62").Value = -3358.5
$ws.Range("N62").Value = -5173

$ws.Range("H65").Value = 3953.75
$ws.Range("I65").Value = 3982.5
$ws.Range("J65").Value = 3925
$ws.Range("K65").Value = 19912.5
$ws.Range("L65").Value = 19625
$ws.Range("M65").Value = -16792.5
$ws.Range("N65").Value = -25865

$ws.Range("H81").Value = 11768274
$ws.Range("I81").Value = 3232
$ws.Range("J81").Value = 16670375
$ws.Range("K81").Value = 6464
$ws.Range("L81").Value = 33340750
$ws.Range("M81").Value = -5403
$ws.Range("N81").Value = -33342872

$ws.Range("H84").Value = 11768274
$ws.Range("I84").Value = 3232
$ws.Range("J84").Value = 16670375
$ws.Range("K84").Value = 32320
$ws.Range("L84").Value = 166703750
$ws.Range("M84").Value = -27016
$ws.Range("N84").Value = -166714358

$ws.Range("H132").Value = 2946.4546
$ws.Range("I132").Value = 2758.8948
$ws.Range("J132").Value = 4134.3335
$ws.Range("K132").Value = 8276.6844
$ws.Range("L132").Value = 12403.0005
$ws.Range("M132").Value = -5746.6844
$ws.Range("N132").Value = -17463.0005
